$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B32").Value = 206822
$ws.Range("C32").Value = -532347
$ws.Range("E32").Value = 535493
$ws.Range("J32").Value = 739184
$ws.Range("P32").Value = -160685
$ws.Range("T32").Value = -160685
$ws.Range("V32").Value = 172962
$ws.Range("X32").Value = 367507
